# fly_record.xlsx update — "badLOC detection and outlier message"
#
# Appends fly 28 (block 7-10) and fly 29 (block 1-15) trial rows to the
# "Sequential Effects ERP" log sheet, introducing the new LIT-condition
# comment labels ("baseline, regular" / "jittering" / "regular" /
# "alternating to calibrate") used to flag bad-LOC / outlier trials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 184
    $ws.Cells.Item(184, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(184, 1).Value2 = 44594
    $ws.Cells.Item(184, 2).Value2 = 28
    $ws.Cells.Item(184, 3).Value2 = 7
    $ws.Cells.Item(184, 4).Value2 = "LIT"
    $ws.Cells.Item(184, 5).Value2 = 0.005
    $ws.Cells.Item(184, 6).Value2 = 0.035
    $ws.Cells.Item(184, 7).Value2 = 25
    $ws.Cells.Item(184, 8).Value2 = 0
    $ws.Cells.Item(184, 9).Value2 = "baseline, regular"

    # Row 185
    $ws.Cells.Item(185, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(185, 1).Value2 = 44594
    $ws.Cells.Item(185, 2).Value2 = 28
    $ws.Cells.Item(185, 3).Value2 = 8
    $ws.Cells.Item(185, 4).Value2 = "LIT"
    $ws.Cells.Item(185, 5).Value2 = 0.005
    $ws.Cells.Item(185, 6).Value2 = 0.035
    $ws.Cells.Item(185, 7).Value2 = 25
    $ws.Cells.Item(185, 8).Value2 = 0
    $ws.Cells.Item(185, 9).Value2 = "jittering"

    # Row 186
    $ws.Cells.Item(186, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(186, 1).Value2 = 44594
    $ws.Cells.Item(186, 2).Value2 = 28
    $ws.Cells.Item(186, 3).Value2 = 9
    $ws.Cells.Item(186, 4).Value2 = "LIT"
    $ws.Cells.Item(186, 5).Value2 = 0.005
    $ws.Cells.Item(186, 6).Value2 = 0.035
    $ws.Cells.Item(186, 7).Value2 = 25
    $ws.Cells.Item(186, 8).Value2 = 0
    $ws.Cells.Item(186, 9).Value2 = "regular"

    # Row 187
    $ws.Cells.Item(187, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(187, 1).Value2 = 44594
    $ws.Cells.Item(187, 2).Value2 = 28
    $ws.Cells.Item(187, 3).Value2 = 10
    $ws.Cells.Item(187, 4).Value2 = "LIT"
    $ws.Cells.Item(187, 5).Value2 = 0.01
    $ws.Cells.Item(187, 6).Value2 = 0.07
    $ws.Cells.Item(187, 7).Value2 = 12.5
    $ws.Cells.Item(187, 8).Value2 = 0

    # Row 188
    $ws.Cells.Item(188, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(188, 1).Value2 = 44595
    $ws.Cells.Item(188, 2).Value2 = 29
    $ws.Cells.Item(188, 3).Value2 = 1
    $ws.Cells.Item(188, 4).Value2 = "LIT"
    $ws.Cells.Item(188, 5).Formula = "=1*1/200"
    $ws.Cells.Item(188, 6).Formula = "=4*1/200"
    $ws.Cells.Item(188, 7).Formula = "=1/(E188+F188)"
    $ws.Cells.Item(188, 8).Value2 = 0
    $ws.Cells.Item(188, 9).Value2 = "alternating to calibrate"

    # Row 189
    $ws.Cells.Item(189, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(189, 1).Value2 = 44595
    $ws.Cells.Item(189, 2).Value2 = 29
    $ws.Cells.Item(189, 3).Value2 = 2
    $ws.Cells.Item(189, 4).Value2 = "LIT"
    $ws.Cells.Item(189, 5).Value2 = 0.005
    $ws.Cells.Item(189, 6).Value2 = 0.035
    $ws.Cells.Item(189, 7).Value2 = 25
    $ws.Cells.Item(189, 8).Value2 = 0
    $ws.Cells.Item(189, 9).Value2 = "alternating to calibrate"

    # Row 190
    $ws.Cells.Item(190, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(190, 1).Value2 = 44595
    $ws.Cells.Item(190, 2).Value2 = 29
    $ws.Cells.Item(190, 3).Value2 = 3
    $ws.Cells.Item(190, 4).Value2 = "LIT"
    $ws.Cells.Item(190, 5).Value2 = 0.01
    $ws.Cells.Item(190, 6).Value2 = 0.07
    $ws.Cells.Item(190, 7).Value2 = 12.5
    $ws.Cells.Item(190, 8).Value2 = 0
    $ws.Cells.Item(190, 9).Value2 = "alternating to calibrate"

    # Row 191
    $ws.Cells.Item(191, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(191, 1).Value2 = 44595
    $ws.Cells.Item(191, 2).Value2 = 29
    $ws.Cells.Item(191, 3).Value2 = 4
    $ws.Cells.Item(191, 4).Value2 = "LIT"
    $ws.Cells.Item(191, 5).Value2 = 0.005
    $ws.Cells.Item(191, 6).Value2 = 0.035
    $ws.Cells.Item(191, 7).Value2 = 25
    $ws.Cells.Item(191, 8).Value2 = 0
    $ws.Cells.Item(191, 9).Value2 = "regular"

    # Row 192
    $ws.Cells.Item(192, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(192, 1).Value2 = 44595
    $ws.Cells.Item(192, 2).Value2 = 29
    $ws.Cells.Item(192, 3).Value2 = 5
    $ws.Cells.Item(192, 4).Value2 = "LIT"
    $ws.Cells.Item(192, 5).Value2 = 0.01
    $ws.Cells.Item(192, 6).Value2 = 0.07
    $ws.Cells.Item(192, 7).Value2 = 12.5
    $ws.Cells.Item(192, 8).Value2 = 0
    $ws.Cells.Item(192, 9).Value2 = "regular"

    # Row 193
    $ws.Cells.Item(193, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(193, 1).Value2 = 44595
    $ws.Cells.Item(193, 2).Value2 = 29
    $ws.Cells.Item(193, 3).Value2 = 6
    $ws.Cells.Item(193, 4).Value2 = "LIT"
    $ws.Cells.Item(193, 5).Value2 = 0.005
    $ws.Cells.Item(193, 6).Value2 = 0.035
    $ws.Cells.Item(193, 7).Value2 = 25
    $ws.Cells.Item(193, 8).Value2 = 0
    $ws.Cells.Item(193, 9).Value2 = "jittering"

    # Row 194
    $ws.Cells.Item(194, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(194, 1).Value2 = 44595
    $ws.Cells.Item(194, 2).Value2 = 29
    $ws.Cells.Item(194, 3).Value2 = 7
    $ws.Cells.Item(194, 4).Value2 = "LIT"
    $ws.Cells.Item(194, 5).Value2 = 0.01
    $ws.Cells.Item(194, 6).Value2 = 0.07
    $ws.Cells.Item(194, 7).Value2 = 12.5
    $ws.Cells.Item(194, 8).Value2 = 0
    $ws.Cells.Item(194, 9).Value2 = "jittering"

    # Row 195
    $ws.Cells.Item(195, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(195, 1).Value2 = 44595
    $ws.Cells.Item(195, 2).Value2 = 29
    $ws.Cells.Item(195, 3).Value2 = 8
    $ws.Cells.Item(195, 4).Value2 = "LIT"
    $ws.Cells.Item(195, 5).Value2 = 0.005
    $ws.Cells.Item(195, 6).Value2 = 0.035
    $ws.Cells.Item(195, 7).Value2 = 25
    $ws.Cells.Item(195, 8).Value2 = 0
    $ws.Cells.Item(195, 9).Value2 = "Red light on"

    # Row 196
    $ws.Cells.Item(196, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(196, 1).Value2 = 44595
    $ws.Cells.Item(196, 2).Value2 = 29
    $ws.Cells.Item(196, 3).Value2 = 9
    $ws.Cells.Item(196, 4).Value2 = "LIT"
    $ws.Cells.Item(196, 5).Value2 = 0.01
    $ws.Cells.Item(196, 6).Value2 = 0.07
    $ws.Cells.Item(196, 7).Value2 = 12.5
    $ws.Cells.Item(196, 8).Value2 = 0
    $ws.Cells.Item(196, 9).Value2 = "Red light on"

    # Row 197
    $ws.Cells.Item(197, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(197, 1).Value2 = 44595
    $ws.Cells.Item(197, 2).Value2 = 29
    $ws.Cells.Item(197, 3).Value2 = 12
    $ws.Cells.Item(197, 4).Value2 = "LIT"
    $ws.Cells.Item(197, 5).Value2 = 0.01
    $ws.Cells.Item(197, 6).Value2 = 0.07
    $ws.Cells.Item(197, 7).Value2 = 12.5
    $ws.Cells.Item(197, 8).Value2 = 0
    $ws.Cells.Item(197, 9).Value2 = "regular"

    # Row 198
    $ws.Cells.Item(198, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(198, 1).Value2 = 44595
    $ws.Cells.Item(198, 2).Value2 = 29
    $ws.Cells.Item(198, 3).Value2 = 13
    $ws.Cells.Item(198, 4).Value2 = "LIT"
    $ws.Cells.Item(198, 5).Value2 = 0.005
    $ws.Cells.Item(198, 6).Value2 = 0.035
    $ws.Cells.Item(198, 7).Value2 = 25
    $ws.Cells.Item(198, 8).Value2 = 0
    $ws.Cells.Item(198, 9).Value2 = "regular"

    # Row 199
    $ws.Cells.Item(199, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(199, 1).Value2 = 44595
    $ws.Cells.Item(199, 2).Value2 = 29
    $ws.Cells.Item(199, 3).Value2 = 14
    $ws.Cells.Item(199, 4).Value2 = "LIT"
    $ws.Cells.Item(199, 5).Value2 = 0.1
    $ws.Cells.Item(199, 6).Value2 = 0.7
    $ws.Cells.Item(199, 7).Value2 = 1.25
    $ws.Cells.Item(199, 8).Value2 = 0
    $ws.Cells.Item(199, 9).Value2 = "jittering"

    # Row 200
    $ws.Cells.Item(200, 1).NumberFormat = "d/mm/yy;@"
    $ws.Cells.Item(200, 1).Value2 = 44595
    $ws.Cells.Item(200, 2).Value2 = 29
    $ws.Cells.Item(200, 3).Value2 = 15
    $ws.Cells.Item(200, 4).Value2 = "LIT"
    $ws.Cells.Item(200, 5).Value2 = 0.1
    $ws.Cells.Item(200, 6).Value2 = 0.7
    $ws.Cells.Item(200, 7).Value2 = 1.25
    $ws.Cells.Item(200, 8).Value2 = 0
    $ws.Cells.Item(200, 9).Value2 = "baseline"

# Restore the book-keeping view state (scroll position + active cell)
# that Excel persists in sheetView/selection after the edits.
$win = $excel.ActiveWindow
$win.ScrollRow = 182
$win.ScrollColumn = 1
$ws.Range("F197").Select()
